$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 763, shifting existing rows (763-804) down to (764-805)
$ws.Rows.Item(763).Insert()

# Populate the newly inserted row with the new log entry.
# Column A holds a date-like string ("2026/02/03") that must stay literal text
# (matching the rest of the sheet), so force Text format while assigning it,
# then clear the format again so the cell carries no explicit style - exactly
# like its neighbouring cells.
$ws.Cells.Item(763, 1).NumberFormat = "@"
$ws.Cells.Item(763, 1).Value = "2026/02/03"
$ws.Cells.Item(763, 1).ClearFormats()

$ws.Cells.Item(763, 2).Value = "火"
$ws.Cells.Item(763, 3).Value = 13
$ws.Cells.Item(763, 4).Value = 201
